# Apply scheduled-runner market data refresh to Lamia_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row from source workbook (original context row 2601)
$ws.Range("H40").Value = 5240.8887
$ws.Range("I40").Value = 3239.6
$ws.Range("J40").Value = 7742.5
$ws.Range("K40").Value = 3239.6
$ws.Range("L40").Value = 7742.5
$ws.Range("M40").Value = -3064.6
$ws.Range("N40").Value = -8092.5

# Row from source workbook (original context row 4116)
$ws.Range("H70").Value = 4178.9
$ws.Range("I70").Value = 1166
$ws.Range("J70").Value = 5470.143
$ws.Range("K70").Value = 3498
$ws.Range("L70").Value = 16410.429
$ws.Range("M70").Value = -3228
$ws.Range("N70").Value = -16950.429

# Row from source workbook (original context row 4269)
$ws.Range("H73").Value = 4178.9
$ws.Range("I73").Value = 1166
$ws.Range("J73").Value = 5470.143
$ws.Range("K73").Value = 3498
$ws.Range("L73").Value = 16410.429
$ws.Range("M73").Value = -2562
$ws.Range("N73").Value = -18282.429

# Row from source workbook (original context row 4422)
$ws.Range("H76").Value = 17666.666

# Row from source workbook (original context row 4575)
$ws.Range("H79").Value = 17666.666

# Row from source workbook (original context row 6427)
$ws.Range("H116").Value = 10464.667
$ws.Range("I116").Value = 9621.666999999999
$ws.Range("J116").Value = 11307.667
$ws.Range("K116").Value = 9621.666999999999
$ws.Range("L116").Value = 11307.667
$ws.Range("M116").Value = -6179.666999999999
$ws.Range("N116").Value = -18191.667

# Row from source workbook (original context row 7223)
$ws.Range("H132").Value = 1574.3125
$ws.Range("I132").Value = 1590.5814
$ws.Range("J132").Value = 1434.4
$ws.Range("K132").Value = 4771.7442
$ws.Range("L132").Value = 4303.200000000001
$ws.Range("M132").Value = -2241.7442
$ws.Range("N132").Value = -9363.200000000001

# Row from source workbook (original context row 7474)
$ws.Range("H137").Value = 31253416
$ws.Range("I137").Value = 76925720
$ws.Range("J137").Value = 3943.7896
$ws.Range("K137").Value = 230777160
$ws.Range("L137").Value = 11831.3688
$ws.Range("M137").Value = -230774610
$ws.Range("N137").Value = -16931.3688

# Row from source workbook (original context row 7673)
$ws.Range("H141").Value = 846.6
$ws.Range("I141").Value = 846.6
$ws.Range("K141").Value = 2539.8
$ws.Range("M141").Value = 2640.2


$ws = $wb.Worksheets.Item("ARM")
# Row from source workbook (original context row 7816)
$ws.Range("H2").Value = 8908.286
$ws.Range("I2").Value = 1250
$ws.Range("K2").Value = 1250
$ws.Range("M2").Value = -1137

# Row from source workbook (original context row 7969)
$ws.Range("H5").Value = 206.92857
$ws.Range("I5").Value = 214.71428
$ws.Range("K5").Value = 214.71428
$ws.Range("M5").Value = -102.71428

# Row from source workbook (original context row 9307)
$ws.Range("H32").Value = 25645012
$ws.Range("I32").Value = 27030364
$ws.Range("J32").Value = 16007
$ws.Range("K32").Value = 27030364
$ws.Range("L32").Value = 16007
$ws.Range("M32").Value = -27030077
$ws.Range("N32").Value = -16581

# Row from source workbook (original context row 9959)
$ws.Range("H45").Value = 4242.385
$ws.Range("I45").Value = 2592.75
$ws.Range("K45").Value = 2592.75
$ws.Range("M45").Value = -2215.75

# Row from source workbook (original context row 13393)
$ws.Range("H116").Value = 8908.286
$ws.Range("I116").Value = 1250
$ws.Range("K116").Value = 1250
$ws.Range("M116").Value = 1044


$ws = $wb.Worksheets.Item("BSM")
# Row from source workbook (original context row 14819)
$ws.Range("H3").Value = 8908.286
$ws.Range("I3").Value = 1250
$ws.Range("K3").Value = 1250
$ws.Range("M3").Value = -1136

# Row from source workbook (original context row 14871)
$ws.Range("H4").Value = 206.92857
$ws.Range("I4").Value = 214.71428
$ws.Range("K4").Value = 214.71428
$ws.Range("M4").Value = -99.71428

# Row from source workbook (original context row 21150)
$ws.Range("H133").Value = 89775
$ws.Range("J133").Value = 89775
$ws.Range("L133").Value = 89775
$ws.Range("N133").Value = -99895

# Row from source workbook (original context row 21199)
$ws.Range("H134").Value = 3116.3
$ws.Range("I134").Value = 1794.3334
$ws.Range("K134").Value = 5383.0002
$ws.Range("M134").Value = -2848.0002

# Row from source workbook (original context row 21346)
$ws.Range("H137").Value = 69991.336
$ws.Range("J137").Value = 69991.336
$ws.Range("L137").Value = 69991.336
$ws.Range("N137").Value = -80191.336

# Row from source workbook (original context row 21395)
$ws.Range("H138").Value = 65491.2
$ws.Range("J138").Value = 65491.2
$ws.Range("L138").Value = 65491.2
$ws.Range("N138").Value = -75771.2


$ws = $wb.Worksheets.Item("CRP")
# Row from source workbook (original context row 24447)
$ws.Range("H58").Value = 5459.8823
$ws.Range("I58").Value = 2177.2222
$ws.Range("K58").Value = 2177.2222
$ws.Range("M58").Value = -1974.2222

# Row from source workbook (original context row 25798)
$ws.Range("H86").Value = 11937.875
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 12786.143
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 12786.143
$ws.Range("M86").Value = -4877
$ws.Range("N86").Value = -15032.143

# Row from source workbook (original context row 25945)
$ws.Range("H89").Value = 11937.875
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 12786.143
$ws.Range("K89").Value = 30000
$ws.Range("L89").Value = 63930.715
$ws.Range("M89").Value = -24384
$ws.Range("N89").Value = -75162.715

# Row from source workbook (original context row 26836)
$ws.Range("H107").Value = 1343.7
$ws.Range("I107").Value = 1362.5714
$ws.Range("J107").Value = 1299.6666
$ws.Range("K107").Value = 1362.5714
$ws.Range("L107").Value = 1299.6666
$ws.Range("M107").Value = 557.4286
$ws.Range("N107").Value = -5139.6666

# Row from source workbook (original context row 28085)
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# Row from source workbook (original context row 28189)
$ws.Range("H134").Value = 9442.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 9442.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 28327.5
$ws.Range("N134").Value = -33397.5
$ws.Range("M134").ClearContents()

# Row from source workbook (original context row 28290)
$ws.Range("H136").Value = 5459.8823
$ws.Range("I136").Value = 2177.2222
$ws.Range("K136").Value = 6531.6666
$ws.Range("M136").Value = -3981.6666


$ws = $wb.Worksheets.Item("CUL")
# Row from source workbook (original context row 28678)
$ws.Range("H2").Value = 220.25
$ws.Range("I2").Value = 168.11111
$ws.Range("J2").Value = 376.66666
$ws.Range("K2").Value = 1008.66666
$ws.Range("L2").Value = 2259.99996
$ws.Range("M2").Value = -895.66666
$ws.Range("N2").Value = -2485.99996

# Row from source workbook (original context row 29737)
$ws.Range("H23").Value = 1055.2
$ws.Range("J23").Value = 1475.1666
$ws.Range("L23").Value = 4425.4998
$ws.Range("N23").Value = -4895.4998

# Row from source workbook (original context row 30499)
$ws.Range("H38").Value = 42.285713
$ws.Range("I38").Value = 52.142857
$ws.Range("J38").Value = 32.42857
$ws.Range("K38").Value = 156.428571
$ws.Range("L38").Value = 97.28570999999999
$ws.Range("M38").Value = 190.571429
$ws.Range("N38").Value = -791.28571


$ws = $wb.Worksheets.Item("GSM")
# Row from source workbook (original context row 39204)
$ws.Range("H70").Value = 10899.608
$ws.Range("I70").Value = 8051.3076
$ws.Range("J70").Value = 14602.4
$ws.Range("K70").Value = 8051.3076
$ws.Range("L70").Value = 14602.4
$ws.Range("M70").Value = -7781.3076
$ws.Range("N70").Value = -15142.4

# Row from source workbook (original context row 39357)
$ws.Range("H73").Value = 10899.608
$ws.Range("I73").Value = 8051.3076
$ws.Range("J73").Value = 14602.4
$ws.Range("K73").Value = 8051.3076
$ws.Range("L73").Value = 14602.4
$ws.Range("M73").Value = -7115.3076
$ws.Range("N73").Value = -16474.4

# Row from source workbook (original context row 41017)
$ws.Range("H107").Value = 1408.2222
$ws.Range("I107").Value = 334.25
$ws.Range("K107").Value = 334.25
$ws.Range("M107").Value = 1585.75

# Row from source workbook (original context row 42230)
$ws.Range("H132").Value = 599833.8
$ws.Range("I132").Value = 693122.9399999999
$ws.Range("K132").Value = 2079368.82
$ws.Range("M132").Value = -2076838.82


$ws = $wb.Worksheets.Item("LTW")
# Row from source workbook (original context row 43503)
$ws.Range("H16").Value = 1856.8572
$ws.Range("I16").Value = 1856.8572
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1856.8572
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1686.8572
$ws.Range("N16").ClearContents()

# Row from source workbook (original context row 46066)
$ws.Range("H68").Value = 7018.875
$ws.Range("I68").Value = 5024.5
$ws.Range("K68").Value = 5024.5
$ws.Range("M68").Value = -4275.5

# Row from source workbook (original context row 46216)
$ws.Range("H71").Value = 7018.875
$ws.Range("I71").Value = 5024.5
$ws.Range("K71").Value = 25122.5
$ws.Range("M71").Value = -21378.5

# Row from source workbook (original context row 49371)
$ws.Range("H136").Value = 6379.2285
$ws.Range("I136").Value = 2302.875
$ws.Range("K136").Value = 6908.625
$ws.Range("M136").Value = -4358.625


$ws = $wb.Worksheets.Item("WVR")
# Row from source workbook (original context row 56108)
$ws.Range("H132").Value = 5572.0977
$ws.Range("I132").Value = 2273.963
$ws.Range("K132").Value = 6821.889000000001
$ws.Range("M132").Value = -4291.889000000001

